# Append one new row (row 55) to each of the 4 worksheets. The new row is
# a copy of the previous last row (row 54) with the timestamp in column A
# advanced by exactly one day - continuing the existing daily log pattern.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $lastRow = 54
    $newRow = 55

    # Column A: timestamp (numeric, date-formatted) - one day after the
    # previous row.
    $prevTime = $ws.Cells.Item($lastRow, 1).Value2
    $ws.Cells.Item($newRow, 1).Value2 = $prevTime + 1
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

    # Columns B..I: copy verbatim from the previous row (text or numeric).
    for ($col = 2; $col -le 9; $col++) {
        $srcCell = $ws.Cells.Item($lastRow, $col)
        $dstCell = $ws.Cells.Item($newRow, $col)
        $dstCell.Value2 = $srcCell.Value2
    }
}
